$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the original "INPSDR0120ITHACABAKERY" / "Cup - Hot (12oz)" line item).
# This shifts all subsequent rows up by one, matching the target layout, and
# collapses the used range from A1:E9 down to A1:E8.
$ws.Rows.Item(2).Delete()
